$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 375, shifting rows 375:480 down to 376:481
$ws.Rows.Item(375).Insert()

# Populate the newly inserted row 375 with its data
$ws.Range("A375").Value = 3
$ws.Range("B375").Value = "Femacal de La Calera"
$ws.Range("C375").Value = "Coquimbo"
$ws.Range("D375").Value = 44722
$ws.Range("E375").Value = 5
$ws.Range("F375").Value = 100112032
$ws.Range("G375").Value = "Zapallo italiano"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 190
$ws.Range("K375").Value = 11000
$ws.Range("L375").Value = 12000
$ws.Range("M375").Value = 11500
$ws.Range("N375").Value = "$/caja 70 unidades"
$ws.Range("O375").Value = "Región de Arica y Parinacota"
$ws.Range("P375").Value = 164
$ws.Range("Q375").Value = 70
$ws.Range("R375").Value = "Hortaliza"

# Match the D-column (date) style used by the rest of the column (style index 2)
$ws.Range("D375").NumberFormat = $ws.Range("D376").NumberFormat
